$wb = $excel.ActiveWorkbook

$timestamp = "2025-12-05 07:07:36"

# Column AA ("as_of_utc") holds a refresh timestamp that needs to be bumped
# for every data row (rows 2-26) on both data sheets: "Главные" and "Линейные".
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 26; $r++) {
        $ws.Cells.Item($r, 27).Value = $timestamp
    }
}
